$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.92947733333333
$ws.Range("H2").Value = 80.788432
$ws.Range("I2").Value = 0.9279949792877585
$ws.Range("J2").Value = 0.9279949792877585
$ws.Range("M2").Value = 91.51130433333333
$ws.Range("N2").Value = 274.533913
$ws.Range("O2").Value = 0.9685519820468944
$ws.Range("P2").Value = 0.9685519820468945
$ws.Range("Q2").Value = 2464.351595788268
$ws.Range("R2").Value = 22179.16436209442
$ws.Range("S2").Value = 0.8988113765187252
$ws.Range("T2").Value = 0.8988113765187253
$ws.Range("G3").Value = 26.92947733333333
$ws.Range("H3").Value = 80.788432
$ws.Range("I3").Value = 0.9279949792877585
$ws.Range("J3").Value = 0.9279949792877585
$ws.Range("O3").Value = 0.001425786415744213
$ws.Range("P3").Value = 0.001425786415744214
$ws.Range("Q3").Value = 3.627723750528
$ws.Range("R3").Value = 32.649513754752
$ws.Range("S3").Value = 0.001323122635347319
$ws.Range("T3").Value = 0.001323122635347319
$ws.Range("G4").Value = 26.92947733333333
$ws.Range("H4").Value = 80.788432
$ws.Range("I4").Value = 0.9279949792877585
$ws.Range("J4").Value = 0.9279949792877585
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509734999999999
$ws.Range("O4").Value = 0.03002223153736139
$ws.Range("P4").Value = 0.03002223153736139
$ws.Range("Q4").Value = 76.38757193172444
$ws.Range("R4").Value = 687.4881473855199
$ws.Range("S4").Value = 0.02786048013368598
$ws.Range("T4").Value = 0.02786048013368598
$ws.Range("I5").Value = 0.04103565698374688
$ws.Range("J5").Value = 0.04103565698374688
$ws.Range("M5").Value = 91.51130433333333
$ws.Range("N5").Value = 274.533913
$ws.Range("O5").Value = 0.9685519820468944
$ws.Range("P5").Value = 0.9685519820468945
$ws.Range("Q5").Value = 108.9728813508578
$ws.Range("R5").Value = 980.7559321577199
$ws.Range("S5").Value = 0.03974516690620453
$ws.Range("T5").Value = 0.03974516690620453
$ws.Range("I6").Value = 0.04103565698374688
$ws.Range("J6").Value = 0.04103565698374688
$ws.Range("O6").Value = 0.001425786415744213
$ws.Range("P6").Value = 0.001425786415744214
$ws.Range("S6").Value = 0.00005850808228856546
$ws.Range("T6").Value = 0.00005850808228856547
$ws.Range("I7").Value = 0.04103565698374688
$ws.Range("J7").Value = 0.04103565698374688
$ws.Range("M7").Value = 2.836578333333333
$ws.Range("N7").Value = 8.509734999999999
$ws.Range("O7").Value = 0.03002223153736139
$ws.Range("P7").Value = 0.03002223153736139
$ws.Range("Q7").Value = 3.377835300377777
$ws.Range("R7").Value = 30.40051770339999
$ws.Range("S7").Value = 0.00123198199525379
$ws.Range("T7").Value = 0.00123198199525379
$ws.Range("G8").Value = 0.8986996666666666
$ws.Range("I8").Value = 0.03096936372849452
$ws.Range("J8").Value = 0.03096936372849452
$ws.Range("M8").Value = 91.51130433333333
$ws.Range("N8").Value = 274.533913
$ws.Range("O8").Value = 0.9685519820468944
$ws.Range("P8").Value = 0.9685519820468945
$ws.Range("Q8").Value = 82.24117870059854
$ws.Range("R8").Value = 740.1706083053868
$ws.Range("S8").Value = 0.02999543862196457
$ws.Range("T8").Value = 0.02999543862196458
$ws.Range("G9").Value = 0.8986996666666666
$ws.Range("I9").Value = 0.03096936372849452
$ws.Range("J9").Value = 0.03096936372849452
$ws.Range("O9").Value = 0.001425786415744213
$ws.Range("P9").Value = 0.001425786415744214
$ws.Range("S9").Value = 0.00004415569810832905
$ws.Range("T9").Value = 0.00004415569810832907
$ws.Range("G10").Value = 0.8986996666666666
$ws.Range("I10").Value = 0.03096936372849452
$ws.Range("J10").Value = 0.03096936372849452
$ws.Range("M10").Value = 2.836578333333333
$ws.Range("N10").Value = 8.509734999999999
$ws.Range("O10").Value = 0.03002223153736139
$ws.Range("P10").Value = 0.03002223153736139
$ws.Range("Q10").Value = 2.549232002640555
$ws.Range("S10").Value = 0.0009297694084216243
$ws.Range("T10").Value = 0.0009297694084216244
